# Auto-generated edit script: apply scheduled-runner value updates to Phantom_Profits sheets
# Updates currentAveragePrice / Leve price & profit columns (H:N) per refreshed market data.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 8231
$ws.Range("I32").Value = 8989
$ws.Range("J32").Value = 7776.2
$ws.Range("K32").Value = 8989
$ws.Range("L32").Value = 7776.2
$ws.Range("M32").Value = -8663
$ws.Range("N32").Value = -8428.200000000001
$ws.Range("H33").Value = 389.2143
$ws.Range("I33").Value = 402.07693
$ws.Range("K33").Value = 402.07693
$ws.Range("M33").Value = -173.07693
$ws.Range("H43").Value = 4575
$ws.Range("I43").Value = 5700.3335
$ws.Range("K43").Value = 5700.3335
$ws.Range("M43").Value = -5631.3335
$ws.Range("H98").Value = 756.75
$ws.Range("I98").Value = 756.75
$ws.Range("K98").Value = 756.75
$ws.Range("M98").Value = 741.25
$ws.Range("H106").Value = 18176.941
$ws.Range("I106").Value = 16334.066
$ws.Range("K106").Value = 16334.066
$ws.Range("M106").Value = -15703.066
$ws.Range("H112").Value = 3968.2856
$ws.Range("J112").Value = 3968.2856
$ws.Range("L112").Value = 11904.8568
$ws.Range("N112").Value = -14120.8568
$ws.Range("H122").Value = 756.75
$ws.Range("I122").Value = 756.75
$ws.Range("K122").Value = 2270.25
$ws.Range("M122").Value = 179.75
$ws.Range("H141").Value = 719
$ws.Range("I141").Value = 821
$ws.Range("J141").Value = 311
$ws.Range("K141").Value = 2463
$ws.Range("L141").Value = 933
$ws.Range("M141").Value = 2717
$ws.Range("N141").Value = -11293

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1282.4
$ws.Range("I74").Value = 978.5
$ws.Range("K74").Value = 978.5
$ws.Range("M74").Value = -104.5
$ws.Range("H77").Value = 1282.4
$ws.Range("I77").Value = 978.5
$ws.Range("K77").Value = 4892.5
$ws.Range("M77").Value = -524.5
$ws.Range("H110").Value = 13149.75
$ws.Range("I110").Value = 13149.75
$ws.Range("K110").Value = 13149.75
$ws.Range("M110").Value = -11104.75
$ws.Range("H122").Value = 1674.36
$ws.Range("I122").Value = 1674.36
$ws.Range("K122").Value = 5023.08
$ws.Range("M122").Value = -2573.08

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -4996

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 1600
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H31").Value = 1789.5
$ws.Range("I31").Value = 1869.5
$ws.Range("J31").Value = 1549.5
$ws.Range("K31").Value = 1869.5
$ws.Range("L31").Value = 1549.5
$ws.Range("M31").Value = -1574.5
$ws.Range("N31").Value = -2139.5
$ws.Range("H34").Value = 1789.5
$ws.Range("I34").Value = 1869.5
$ws.Range("J34").Value = 1549.5
$ws.Range("K34").Value = 1869.5
$ws.Range("L34").Value = 1549.5
$ws.Range("M34").Value = -1667.5
$ws.Range("N34").Value = -1953.5
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H95").Value = 33649.832
$ws.Range("J95").Value = 33649.832
$ws.Range("L95").Value = 33649.832
$ws.Range("N95").Value = -39141.832
$ws.Range("H109").Value = 90999.10000000001
$ws.Range("I109").Value = 90000
$ws.Range("J109").Value = 91110.11
$ws.Range("K109").Value = 90000
$ws.Range("L109").Value = 91110.11
$ws.Range("M109").Value = -88960
$ws.Range("N109").Value = -93190.11

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1469243.2
$ws.Range("I4").Value = 42705.4
$ws.Range("J4").Value = 30000000
$ws.Range("K4").Value = 128116.2
$ws.Range("L4").Value = 90000000
$ws.Range("M4").Value = -128004.2
$ws.Range("N4").Value = -90000224
$ws.Range("H9").Value = 4333.3335
$ws.Range("J9").Value = 4500
$ws.Range("L9").Value = 13500
$ws.Range("N9").Value = -13948
$ws.Range("H12").Value = 356.83334
$ws.Range("I12").Value = 342.75
$ws.Range("K12").Value = 1028.25
$ws.Range("M12").Value = -855.25
$ws.Range("H17").Value = 621.8
$ws.Range("I17").Value = 175
$ws.Range("J17").Value = 919.6667
$ws.Range("K17").Value = 525
$ws.Range("L17").Value = 2759.0001
$ws.Range("M17").Value = -356
$ws.Range("N17").Value = -3097.0001
$ws.Range("H70").Value = 15709
$ws.Range("J70").Value = 16662.666
$ws.Range("L70").Value = 49987.99800000001
$ws.Range("N70").Value = -50617.99800000001
$ws.Range("H73").Value = 15709
$ws.Range("J73").Value = 16662.666
$ws.Range("L73").Value = 49987.99800000001
$ws.Range("N73").Value = -52171.99800000001
$ws.Range("H81").Value = 11673.75
$ws.Range("I81").Value = 1695
$ws.Range("K81").Value = 5085
$ws.Range("M81").Value = -3962
$ws.Range("H84").Value = 11673.75
$ws.Range("I84").Value = 1695
$ws.Range("K84").Value = 15255
$ws.Range("M84").Value = -9639
$ws.Range("H106").Value = 15498.571
$ws.Range("J106").Value = 18500
$ws.Range("L106").Value = 55500
$ws.Range("N106").Value = -57392
$ws.Range("H107").Value = 778.5454999999999
$ws.Range("I107").Value = 781.6
$ws.Range("J107").Value = 776
$ws.Range("K107").Value = 2344.8
$ws.Range("L107").Value = 2328
$ws.Range("M107").Value = -424.8000000000002
$ws.Range("N107").Value = -6168
$ws.Range("H110").Value = 8444
$ws.Range("I110").Value = 8444
$ws.Range("K110").Value = 25332
$ws.Range("M110").Value = -21242
$ws.Range("H112").Value = 7733.875
$ws.Range("I112").Value = 2310.3333
$ws.Range("J112").Value = 10988
$ws.Range("K112").Value = 6930.999899999999
$ws.Range("L112").Value = 32964
$ws.Range("M112").Value = -5822.999899999999
$ws.Range("N112").Value = -35180

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 50000
$ws.Range("I26").Value = 50000
$ws.Range("K26").Value = 50000
$ws.Range("M26").Value = -49720
$ws.Range("H50").Value = 50000
$ws.Range("I50").Value = 50000
$ws.Range("K50").Value = 50000
$ws.Range("M50").Value = -49502
$ws.Range("H92").Value = 12884.6
$ws.Range("J92").Value = 12884.6
$ws.Range("L92").Value = 12884.6
$ws.Range("N92").Value = -16628.6
$ws.Range("H97").Value = 573.8333
$ws.Range("I97").Value = 577.6
$ws.Range("J97").Value = 555
$ws.Range("K97").Value = 577.6
$ws.Range("L97").Value = 555
$ws.Range("M97").Value = -81.60000000000002
$ws.Range("N97").Value = -1547
$ws.Range("H102").Value = 2368.5
$ws.Range("I102").Value = 2368.5
$ws.Range("K102").Value = 2368.5
$ws.Range("M102").Value = -746.5
$ws.Range("H121").Value = 107500
$ws.Range("J121").Value = 107500
$ws.Range("L121").Value = 107500
$ws.Range("N121").Value = -110994
$ws.Range("H122").Value = 2538.476
$ws.Range("I122").Value = 2089.353
$ws.Range("K122").Value = 6268.059
$ws.Range("M122").Value = -3818.059

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H61").Value = 1883.3636
$ws.Range("J61").Value = 1832.3334
$ws.Range("L61").Value = 1832.3334
$ws.Range("N61").Value = -2236.3334
$ws.Range("H93").Value = 1999.75
$ws.Range("I93").Value = 1999.75
$ws.Range("K93").Value = 1999.75
$ws.Range("M93").Value = -751.75
$ws.Range("H113").Value = 1883.3636
$ws.Range("J113").Value = 1832.3334
$ws.Range("L113").Value = 1832.3334
$ws.Range("N113").Value = -6172.3334
$ws.Range("H122").Value = 4042.9473
$ws.Range("I122").Value = 3618.25
$ws.Range("K122").Value = 10854.75
$ws.Range("M122").Value = -8404.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6645.222
$ws.Range("I132").Value = 7961.4
$ws.Range("K132").Value = 23884.2
$ws.Range("M132").Value = -21354.2
$ws.Range("H139").Value = 42000
$ws.Range("J139").Value = 42000
$ws.Range("L139").Value = 42000
$ws.Range("N139").Value = -52280
